# Final touches to data mining
# Applies corrected per-county stats (Total/percentages/engine-CO2-mass averages)
# plus a handful of model-name corrections in the top-3 model columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (HARJU MAAKOND)
$ws.Range("C2").Value = 280727
$ws.Range("D2").Value = 9.56
$ws.Range("E2").Value = 54.11
$ws.Range("F2").Value = 58.44
$ws.Range("G2").Value = 40.87
$ws.Range("H2").Value = 0.35
$ws.Range("K2").Value = 109.8849269931288
$ws.Range("L2").Value = 160.7885123833163
$ws.Range("M2").Value = 1568.788310351338
$ws.Range("N2").Value = "TOYOTA 11.2"
$ws.Range("Q2").Value = "OCTAVIA 2.6"

# Row 3 (TARTU MAAKOND)
$ws.Range("C3").Value = 70275
$ws.Range("D3").Value = 6.95
$ws.Range("E3").Value = 44.49
$ws.Range("F3").Value = 50.49
$ws.Range("G3").Value = 48.86
$ws.Range("H3").Value = 0.21
$ws.Range("K3").Value = 106.2803941657772
$ws.Range("L3").Value = 164.7819673424312
$ws.Range("M3").Value = 1570.171597296336
$ws.Range("N3").Value = "VOLKSWAGEN 15.2"
$ws.Range("P3").Value = "FORD 7.8"

# Row 4 (IDA-VIRU MAAKOND)
$ws.Range("C4").Value = 50557
$ws.Range("E4").Value = 39.11
$ws.Range("F4").Value = 48.23
$ws.Range("G4").Value = 51.57
$ws.Range("H4").Value = 0.05
$ws.Range("K4").Value = 102.8087544751469
$ws.Range("L4").Value = 170.8912771588959
$ws.Range("M4").Value = 1573.435805130842

# Row 5 (PÄRNU MAAKOND)
$ws.Range("C5").Value = 40973
$ws.Range("D5").Value = 7.87
$ws.Range("E5").Value = 43.36
$ws.Range("F5").Value = 51.84
$ws.Range("G5").Value = 47.63
$ws.Range("H5").Value = 0.16
$ws.Range("K5").Value = 104.1971176140385
$ws.Range("L5").Value = 166.8096053091163
$ws.Range("M5").Value = 1550.287994532985

# Row 6 (LÄÄNE-VIRU MAAKOND)
$ws.Range("C6").Value = 29192
$ws.Range("D6").Value = 8.72
$ws.Range("E6").Value = 38.27
$ws.Range("F6").Value = 50.22
$ws.Range("G6").Value = 49.42
$ws.Range("H6").Value = 0.12
$ws.Range("K6").Value = 101.2792271855303
$ws.Range("L6").Value = 166.6256606432778
$ws.Range("M6").Value = 1542.583858591395

# Row 7 (VILJANDI MAAKOND)
$ws.Range("C7").Value = 24223
$ws.Range("D7").Value = 7.99
$ws.Range("E7").Value = 39.47
$ws.Range("F7").Value = 50.7
$ws.Range("G7").Value = 48.8
$ws.Range("H7").Value = 0.25
$ws.Range("K7").Value = 100.8744953143707
$ws.Range("L7").Value = 164.2860212413285
$ws.Range("M7").Value = 1535.333814969244
$ws.Range("O7").Value = "AUDI 8.1"

# Row 8 (VÕRU MAAKOND)
$ws.Range("C8").Value = 19212
$ws.Range("D8").Value = 8.45
$ws.Range("E8").Value = 33.53
$ws.Range("F8").Value = 42.45
$ws.Range("G8").Value = 57.14
$ws.Range("H8").Value = 0.2
$ws.Range("K8").Value = 99.87674370185303
$ws.Range("L8").Value = 165.4419168687798
$ws.Range("M8").Value = 1538.682125754737
$ws.Range("S8").Value = "A6 AVANT 2.9"

# Row 9 (RAPLA MAAKOND)
$ws.Range("C9").Value = 18255
$ws.Range("D9").Value = 8.44
$ws.Range("E9").Value = 42.97
$ws.Range("F9").Value = 50.31
$ws.Range("G9").Value = 49.42
$ws.Range("H9").Value = 0.12
$ws.Range("K9").Value = 102.8584278279924
$ws.Range("L9").Value = 165.4416983523447
$ws.Range("M9").Value = 1551.253190906601

# Row 10 (SAARE MAAKOND)
$ws.Range("C10").Value = 18046
$ws.Range("D10").Value = 5.62
$ws.Range("E10").Value = 40.11
$ws.Range("F10").Value = 58.01
$ws.Range("G10").Value = 41.58
$ws.Range("H10").Value = 0.3
$ws.Range("K10").Value = 100.8285880527541
$ws.Range("L10").Value = 166.520876056107
$ws.Range("M10").Value = 1523.387675939266
$ws.Range("P10").Value = "AUDI 6.3"

# Row 11 (JÄRVA MAAKOND)
$ws.Range("C11").Value = 15461
$ws.Range("D11").Value = 8.72
$ws.Range("E11").Value = 37.26
$ws.Range("F11").Value = 49.25
$ws.Range("G11").Value = 50.42
$ws.Range("H11").Value = 0.19
$ws.Range("K11").Value = 101.7334907185822
$ws.Range("L11").Value = 166.8556081214272
$ws.Range("M11").Value = 1539.461419054395
$ws.Range("N11").Value = "VOLKSWAGEN 13.4"

# Row 12 (MÄÄRAMATA)
$ws.Range("C12").Value = 14867
$ws.Range("D12").Value = 26.41
$ws.Range("E12").Value = 48.6
$ws.Range("F12").Value = 51.85
$ws.Range("G12").Value = 47.97
$ws.Range("H12").Value = 0.05
$ws.Range("K12").Value = 111.5234277258357
$ws.Range("L12").Value = 179.1894533469712
$ws.Range("M12").Value = 1589.406941548396

# Row 13 (JÕGEVA MAAKOND)
$ws.Range("C13").Value = 14454
$ws.Range("D13").Value = 7.86
$ws.Range("E13").Value = 35.59
$ws.Range("F13").Value = 46.02
$ws.Range("G13").Value = 53.58
$ws.Range("H13").Value = 0.17
$ws.Range("K13").Value = 100.2577210460772
$ws.Range("L13").Value = 166.9791041577368
$ws.Range("M13").Value = 1541.765116922651
$ws.Range("N13").Value = "VOLKSWAGEN 15.9"

# Row 14 (PÕLVA MAAKOND)
$ws.Range("C14").Value = 14305
$ws.Range("D14").Value = 9.05
$ws.Range("E14").Value = 32.75
$ws.Range("F14").Value = 45.7
$ws.Range("G14").Value = 53.62
$ws.Range("H14").Value = 0.38
$ws.Range("K14").Value = 99.21154141908424
$ws.Range("L14").Value = 167.1736182925428
$ws.Range("M14").Value = 1537.96798322265
$ws.Range("N14").Value = "VOLKSWAGEN 19.0"
$ws.Range("O14").Value = "AUDI 9.8"
$ws.Range("P14").Value = "FORD 7.8"

# Row 15 (VALGA MAAKOND)
$ws.Range("C15").Value = 14186
$ws.Range("D15").Value = 9.4
$ws.Range("E15").Value = 33.03
$ws.Range("F15").Value = 43.61
$ws.Range("G15").Value = 56.03
$ws.Range("H15").Value = 0.2
$ws.Range("K15").Value = 100.2669956294939
$ws.Range("L15").Value = 168.1528690807799
$ws.Range("M15").Value = 1549.174185817003

# Row 16 (LÄÄNE MAAKOND)
$ws.Range("C16").Value = 11154
$ws.Range("D16").Value = 7.26
$ws.Range("E16").Value = 42.74
$ws.Range("F16").Value = 55.86
$ws.Range("G16").Value = 43.22
$ws.Range("H16").Value = 0.2
$ws.Range("K16").Value = 101.1230141653219
$ws.Range("L16").Value = 162.4164338328677
$ws.Range("M16").Value = 1524.597722790031
$ws.Range("N16").Value = "VOLKSWAGEN 11.0"
$ws.Range("P16").Value = "VOLVO 6.9"

# Row 17 (HIIU MAAKOND)
$ws.Range("C17").Value = 5635
$ws.Range("D17").Value = 5.91
$ws.Range("E17").Value = 40.37
$ws.Range("F17").Value = 54.04
$ws.Range("G17").Value = 45.77
$ws.Range("H17").Value = 0.16
$ws.Range("K17").Value = 101.1832298136646
$ws.Range("L17").Value = 168.9677163270707
$ws.Range("M17").Value = 1540.292990239574
$ws.Range("N17").Value = "VOLKSWAGEN 12.8"
$ws.Range("Q17").Value = "AVENSIS 2.4"
$ws.Range("R17").Value = "PASSAT VARIANT 2.4"

